$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.188.27"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "3.519.47"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "585.88"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "132.39"
$ws.Range("E6").Value = "  -3.76%  "
$ws.Range("D7").Value = "3.521.64"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "4.110.56"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("E15").Value = "  -3.42%  "
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "3.509.70"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "64.202.04"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").Value = "9.99"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("D22").Value = "392.11"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "0.579"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").Value = "3.654.91"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "73.00"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "0.0000112"
$ws.Range("E27").Value = "  -4.77%  "
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").Value = "7.50"
$ws.Range("E29").Value = "  -7.74%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D33").Value = "3.514.48"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "23.96"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").Value = "5.40"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "1.57"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").Value = "168.07"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D42").Value = "26.64"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "41.96"
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").Value = "1.21"
$ws.Range("E46").Value = "  -5.61%  "
$ws.Range("D47").Value = "4.39"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("D49").Value = "2.451.29"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").Value = "0.897"
$ws.Range("E51").Value = "  -1.36%  "
